# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet, per the commit's symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.46"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-0.04%"
$ws.Range("E2").ClearFormats()

$ws.Range("D3").Value = "'44.30"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'1.43%"
$ws.Range("E3").ClearFormats()

$ws.Range("D4").Value = "'5.567"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'2.36%"
$ws.Range("E4").ClearFormats()

$ws.Range("D5").Value = "'0.08064"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-1.07%"
$ws.Range("E5").ClearFormats()

$ws.Range("D6").Value = "'1.913"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-0.02%"
$ws.Range("E6").ClearFormats()

$ws.Range("E7").Value = "'-7.62%"
$ws.Range("E7").ClearFormats()

$ws.Range("D8").Value = "'0.9521"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'1.03%"
$ws.Range("E8").ClearFormats()

$ws.Range("D9").Value = "'0.1183"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'0.64%"
$ws.Range("E9").ClearFormats()

$ws.Range("D10").Value = "'0.1851"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'-2.14%"
$ws.Range("E10").ClearFormats()

$ws.Range("D11").Value = "'10.24"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'16.81%"
$ws.Range("E11").ClearFormats()

$ws.Range("D12").Value = "'0.09761"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'0.08%"
$ws.Range("E12").ClearFormats()

$ws.Range("D13").Value = "'0.04617"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'8.65%"
$ws.Range("E13").ClearFormats()

$ws.Range("D14").Value = "'0.1066"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-0.18%"
$ws.Range("E14").ClearFormats()

$ws.Range("D15").Value = "'0.001287"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'0.63%"
$ws.Range("E15").ClearFormats()

$ws.Range("D16").Value = "'0.04201"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'-3.86%"
$ws.Range("E16").ClearFormats()

$ws.Range("D17").Value = "'0.005966"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'-1.26%"
$ws.Range("E17").ClearFormats()

$ws.Range("D18").Value = "'3.395"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-4.28%"
$ws.Range("E18").ClearFormats()

$ws.Range("D19").Value = "'4.301"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'-0.35%"
$ws.Range("E19").ClearFormats()

$ws.Range("E20").Value = "'-1.55%"
$ws.Range("E20").ClearFormats()

$ws.Range("E21").Value = "'4.41%"
$ws.Range("E21").ClearFormats()

$ws.Range("D22").Value = "'0.2507"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'0.24%"
$ws.Range("E22").ClearFormats()

$ws.Range("D23").Value = "'0.001247"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'0.61%"
$ws.Range("E23").ClearFormats()

$ws.Range("D24").Value = "'0.004354"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'0.09%"
$ws.Range("E24").ClearFormats()

$ws.Range("E25").Value = "'-3.63%"
$ws.Range("E25").ClearFormats()

$ws.Range("E26").Value = "'-0.70%"
$ws.Range("E26").ClearFormats()

$ws.Range("D38").Value = "'0.02683"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'0.30%"
$ws.Range("E38").ClearFormats()

$ws.Range("D39").Value = "'0.05540"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'-1.74%"
$ws.Range("E39").ClearFormats()

$ws.Range("D40").Value = "'0.007566"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'-4.02%"
$ws.Range("E40").ClearFormats()

$ws.Range("D41").Value = "'0.1405"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-0.11%"
$ws.Range("E41").ClearFormats()

$ws.Range("D42").Value = "'0.008329"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-14.81%"
$ws.Range("E42").ClearFormats()

$ws.Range("E43").Value = "'-4.90%"
$ws.Range("E43").ClearFormats()

$ws.Range("D44").Value = "'0.008911"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-7.38%"
$ws.Range("E44").ClearFormats()

$ws.Range("D45").Value = "'0.00007108"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'1.01%"
$ws.Range("E45").ClearFormats()

$ws.Range("E46").Value = "'-0.41%"
$ws.Range("E46").ClearFormats()

$ws.Range("D47").Value = "'0.002995"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'-13.59%"
$ws.Range("E47").ClearFormats()

$ws.Range("D48").Value = "'0.002272"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'-0.38%"
$ws.Range("E48").ClearFormats()

$ws.Range("E49").Value = "'-0.41%"
$ws.Range("E49").ClearFormats()

$ws.Range("E50").Value = "'-0.41%"
$ws.Range("E50").ClearFormats()
